$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D2").Value = 11.93
$ws.Range("F2").Value = 7

$ws.Range("E3").Value = 10.16
$ws.Range("F3").Value = 9.949999999999999

$ws.Range("B4").Value = 8.07
$ws.Range("E4").Value = 9.949999999999999

$ws.Range("C5").Value = 9.84
$ws.Range("D5").Value = 10.05
$ws.Range("F5").Value = 10.32

$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 10.05
$ws.Range("E6").Value = 9.68
$ws.Range("G6").Value = 10.23
$ws.Range("H6").Value = 11.63

$ws.Range("F7").Value = 9.77

$ws.Range("F8").Value = 8.369999999999999
